$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 645 (this shifts the existing rows 645:679
# down to 648:682 and grows the used range accordingly).
$ws.Rows("645:647").Insert()

# Populate the 3 newly inserted rows with the new weekly Kiwi price records
# (Femacal de La Calera, Coquimbo, Hayward, Especial/Primera/Segunda).

# Row 645 - Especial
$ws.Cells.Item(645, 1).Value = 3
$ws.Cells.Item(645, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(645, 3).Value = "Coquimbo"
$ws.Cells.Item(645, 4).Value = 44753
$ws.Cells.Item(645, 5).Value = 5
$ws.Cells.Item(645, 6).Value = "Fruta"
$ws.Cells.Item(645, 7).Value = 100101
$ws.Cells.Item(645, 8).Value = "Berries"
$ws.Cells.Item(645, 9).Value = 100101007
$ws.Cells.Item(645, 10).Value = "Kiwi"
$ws.Cells.Item(645, 11).Value = "Hayward"
$ws.Cells.Item(645, 12).Value = "Especial"
$ws.Cells.Item(645, 13).Value = 75
$ws.Cells.Item(645, 14).Value = 7000
$ws.Cells.Item(645, 15).Value = 7000
$ws.Cells.Item(645, 16).Value = 7000
$ws.Cells.Item(645, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(645, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(645, 19).Value = 700
$ws.Cells.Item(645, 20).Value = 10

# Row 646 - Primera
$ws.Cells.Item(646, 1).Value = 3
$ws.Cells.Item(646, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(646, 3).Value = "Coquimbo"
$ws.Cells.Item(646, 4).Value = 44753
$ws.Cells.Item(646, 5).Value = 5
$ws.Cells.Item(646, 6).Value = "Fruta"
$ws.Cells.Item(646, 7).Value = 100101
$ws.Cells.Item(646, 8).Value = "Berries"
$ws.Cells.Item(646, 9).Value = 100101007
$ws.Cells.Item(646, 10).Value = "Kiwi"
$ws.Cells.Item(646, 11).Value = "Hayward"
$ws.Cells.Item(646, 12).Value = "Primera"
$ws.Cells.Item(646, 13).Value = 80
$ws.Cells.Item(646, 14).Value = 6000
$ws.Cells.Item(646, 15).Value = 6000
$ws.Cells.Item(646, 16).Value = 6000
$ws.Cells.Item(646, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(646, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(646, 19).Value = 600
$ws.Cells.Item(646, 20).Value = 10

# Row 647 - Segunda
$ws.Cells.Item(647, 1).Value = 3
$ws.Cells.Item(647, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(647, 3).Value = "Coquimbo"
$ws.Cells.Item(647, 4).Value = 44753
$ws.Cells.Item(647, 5).Value = 5
$ws.Cells.Item(647, 6).Value = "Fruta"
$ws.Cells.Item(647, 7).Value = 100101
$ws.Cells.Item(647, 8).Value = "Berries"
$ws.Cells.Item(647, 9).Value = 100101007
$ws.Cells.Item(647, 10).Value = "Kiwi"
$ws.Cells.Item(647, 11).Value = "Hayward"
$ws.Cells.Item(647, 12).Value = "Segunda"
$ws.Cells.Item(647, 13).Value = 70
$ws.Cells.Item(647, 14).Value = 5000
$ws.Cells.Item(647, 15).Value = 5000
$ws.Cells.Item(647, 16).Value = 5000
$ws.Cells.Item(647, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(647, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(647, 19).Value = 500
$ws.Cells.Item(647, 20).Value = 10
